# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data updates to the Seraph_Profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 198.84616
$ws.Range("I33").Value = 165.41667
$ws.Range("K33").Value = 165.41667
$ws.Range("M33").Value = 63.58332999999999

$ws.Range("H34").Value = 3378.7
$ws.Range("I34").Value = 3378.7
$ws.Range("K34").Value = 3378.7
$ws.Range("M34").Value = -3175.7

$ws.Range("H36").Value = 3378.7
$ws.Range("I36").Value = 3378.7
$ws.Range("K36").Value = 3378.7
$ws.Range("M36").Value = -2663.7

$ws.Range("H138").Value = 2774.8696
$ws.Range("J138").Value = 3726.5454
$ws.Range("L138").Value = 11179.6362
$ws.Range("N138").Value = -21459.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1306.6666
$ws.Range("I5").Value = 960
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 960
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = -848
$ws.Range("N5").Value = -2224

$ws.Range("H102").Value = 18527586
$ws.Range("I102").Value = 22232702
$ws.Range("K102").Value = 22232702
$ws.Range("M102").Value = -22231080

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1306.6666
$ws.Range("I4").Value = 960
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 960
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = -845
$ws.Range("N4").Value = -2230

$ws.Range("H22").Value = 148.2
$ws.Range("I22").Value = 60.25
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 60.25
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = 112.75
$ws.Range("N22").Value = -846

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 89426.336
$ws.Range("I22").Value = 149424.88
$ws.Range("K22").Value = 149424.88
$ws.Range("M22").Value = -149074.88

$ws.Range("H31").Value = 3716.9333
$ws.Range("I31").Value = 2198.6667
$ws.Range("J31").Value = 5994.3335
$ws.Range("K31").Value = 2198.6667
$ws.Range("L31").Value = 5994.3335
$ws.Range("M31").Value = -1903.6667
$ws.Range("N31").Value = -6584.3335

$ws.Range("H34").Value = 3716.9333
$ws.Range("I34").Value = 2198.6667
$ws.Range("J34").Value = 5994.3335
$ws.Range("K34").Value = 2198.6667
$ws.Range("L34").Value = 5994.3335
$ws.Range("M34").Value = -1996.6667
$ws.Range("N34").Value = -6398.3335

$ws.Range("H58").Value = 2506.1875
$ws.Range("I58").Value = 1734.9231
$ws.Range("K58").Value = 1734.9231
$ws.Range("M58").Value = -1531.9231

$ws.Range("H132").Value = 1788.8
$ws.Range("I132").Value = 1765.4445
$ws.Range("J132").Value = 1999
$ws.Range("K132").Value = 5296.333500000001
$ws.Range("L132").Value = 5997
$ws.Range("M132").Value = -2766.333500000001
$ws.Range("N132").Value = -11057

$ws.Range("H136").Value = 2506.1875
$ws.Range("I136").Value = 1734.9231
$ws.Range("K136").Value = 5204.7693
$ws.Range("M136").Value = -2654.7693

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 286087.72
$ws.Range("I29").Value = 666701.7
$ws.Range("K29").Value = 2000105.1
$ws.Range("M29").Value = -1999828.1

$ws.Range("H39").Value = 11573.538
$ws.Range("J39").Value = 11573.538
$ws.Range("L39").Value = 34720.614
$ws.Range("N39").Value = -35308.614

$ws.Range("H116").Value = 3000
$ws.Range("J116").Value = 3000
$ws.Range("L116").Value = 9000
$ws.Range("N116").Value = -15884

$ws.Range("H118").Value = 720
$ws.Range("I118").Value = 720
$ws.Range("K118").Value = 2160
$ws.Range("M118").Value = -917

$ws.Range("H121").Value = 945.7692
$ws.Range("I121").Value = 332
$ws.Range("J121").Value = 1129.9
$ws.Range("K121").Value = 996
$ws.Range("L121").Value = 3389.7
$ws.Range("M121").Value = 314
$ws.Range("N121").Value = -6009.700000000001

$ws.Range("H129").Value = 938.7
$ws.Range("I129").Value = 629.1429000000001
$ws.Range("K129").Value = 1887.4287
$ws.Range("M129").Value = 3112.5713

$ws.Range("H136").Value = 7523
$ws.Range("I136").Value = 1284.5
$ws.Range("K136").Value = 3853.5
$ws.Range("M136").Value = 1246.5

$ws.Range("H138").Value = 5374.8335
$ws.Range("I138").Value = 3832.6667
$ws.Range("J138").Value = 5888.8887
$ws.Range("K138").Value = 11498.0001
$ws.Range("L138").Value = 17666.6661
$ws.Range("M138").Value = -6358.000100000001
$ws.Range("N138").Value = -27946.6661

$ws.Range("H139").Value = 1384.25
$ws.Range("I139").Value = 1216.1818
$ws.Range("J139").Value = 3233
$ws.Range("K139").Value = 3648.5454
$ws.Range("L139").Value = 9699
$ws.Range("M139").Value = 1491.4546
$ws.Range("N139").Value = -19979

$ws.Range("H140").Value = 2775.8333
$ws.Range("I140").Value = 2331.4
$ws.Range("K140").Value = 6994.200000000001
$ws.Range("M140").Value = -1814.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 381.875
$ws.Range("I2").Value = 101.666664
$ws.Range("J2").Value = 550
$ws.Range("K2").Value = 101.666664
$ws.Range("L2").Value = 550
$ws.Range("M2").Value = 11.333336
$ws.Range("N2").Value = -776

$ws.Range("H35").Value = 2508206.2
$ws.Range("J35").Value = 1250555.5
$ws.Range("L35").Value = 1250555.5
$ws.Range("N35").Value = -1251151.5

$ws.Range("H102").Value = 1385.8572
$ws.Range("I102").Value = 1401.1818
$ws.Range("K102").Value = 1401.1818
$ws.Range("M102").Value = 220.8181999999999

$ws.Range("H132").Value = 3138.5667
$ws.Range("I132").Value = 3146.5557
$ws.Range("K132").Value = 9439.667099999999
$ws.Range("M132").Value = -6909.667099999999

$ws.Range("H135").Value = 101060.336
$ws.Range("I135").Value = 101060.336
$ws.Range("K135").Value = 101060.336
$ws.Range("M135").Value = -95990.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7933.6
$ws.Range("I136").Value = 7502
$ws.Range("K136").Value = 22506
$ws.Range("M136").Value = -19956

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 24291.334
$ws.Range("J25").Value = 24291.334
$ws.Range("L25").Value = 24291.334
$ws.Range("N25").Value = -24877.334

$ws.Range("H47").Value = 14641.857
$ws.Range("I47").Value = 9000
$ws.Range("J47").Value = 15582.167
$ws.Range("K47").Value = 9000
$ws.Range("L47").Value = 15582.167
$ws.Range("M47").Value = -8428
$ws.Range("N47").Value = -16726.167

$ws.Range("H54").Value = 22962
$ws.Range("I54").Value = 2000
$ws.Range("K54").Value = 2000
$ws.Range("M54").Value = -1480

$ws.Range("H122").Value = 3962.7144
$ws.Range("I122").Value = 3927.7144
$ws.Range("K122").Value = 11783.1432
$ws.Range("M122").Value = -9333.143199999999

$ws.Range("H132").Value = 926.7143
$ws.Range("I132").Value = 926.7143
$ws.Range("K132").Value = 2780.1429
$ws.Range("M132").Value = -250.1428999999998

$ws.Range("H136").Value = 8741.25
$ws.Range("I136").Value = 9332
$ws.Range("K136").Value = 27996
$ws.Range("M136").Value = -25446

